$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert this week's new observation at the top of the
# Berenjena series (row 471) and push the previous 22 weeks' rows down by
# one (471-492 -> 472-493), extending the used range to row 493.
for ($r = 492; $r -ge 471; $r--) {
    $destRow = $r + 1
    $ws.Range("A" + $destRow + ":R" + $destRow).Value2 = $ws.Range("A" + $r + ":R" + $r).Value2
}

# Row 493 is brand new, so the date cell needs the same number format as the
# rest of column D (copying into an existing row carries this over already).
$ws.Range("D493").NumberFormat = $ws.Range("D492").NumberFormat

# Write this week's new record into row 471.
$ws.Cells.Item(471, 1).Value = 10
$ws.Cells.Item(471, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(471, 3).Value = "La Araucanía"
$ws.Cells.Item(471, 4).Value = 45147
$ws.Cells.Item(471, 5).Value = 9
$ws.Cells.Item(471, 6).Value = 100112001
$ws.Cells.Item(471, 7).Value = "Berenjena"
$ws.Cells.Item(471, 8).Value = "Sin especificar"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 55
$ws.Cells.Item(471, 11).Value = 14000
$ws.Cells.Item(471, 12).Value = 14000
$ws.Cells.Item(471, 13).Value = 14000
$ws.Cells.Item(471, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(471, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(471, 16).Value = 350
$ws.Cells.Item(471, 17).Value = 40
$ws.Cells.Item(471, 18).Value = "Hortaliza"
